$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 56472
$ws.Range("B2").Value = "Antônio Carvalho"
$ws.Range("C2").Value = "Operações"
$ws.Range("D2").Value = "Doença"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45095
$ws.Range("G2").Value = 7049.48

# Row 3
$ws.Range("A3").Value = 84746
$ws.Range("B3").Value = "Ana Luiza da Paz"
$ws.Range("C3").Value = "P&D"
$ws.Range("F3").Value = 45084
$ws.Range("G3").Value = 11758.98

# Row 4
$ws.Range("A4").Value = 70974
$ws.Range("B4").Value = "Maria Vitória Carvalho"
$ws.Range("C4").Value = "Recursos Humanos"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45090
$ws.Range("G4").Value = 3690.48

# Row 5
$ws.Range("A5").Value = 22494
$ws.Range("B5").Value = "Dr. Felipe Sales"
$ws.Range("C5").Value = "Vendas"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45106
$ws.Range("G5").Value = 8039.21

# Row 6
$ws.Range("A6").Value = 20938
$ws.Range("B6").Value = "Vitor da Rosa"
$ws.Range("C6").Value = "Atendimento ao Cliente"
$ws.Range("D6").Value = "Viagem de negócios"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 45097
$ws.Range("G6").Value = 3221.88

# Row 7
$ws.Range("A7").Value = 50526
$ws.Range("B7").Value = "Dr. Luiz Fernando Jesus"
$ws.Range("D7").Value = "Doença"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 45087
$ws.Range("G7").Value = 5457.11

# Row 8
$ws.Range("A8").Value = 25868
$ws.Range("B8").Value = "Dr. Lorenzo Almeida"
$ws.Range("C8").Value = "Financeiro"
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 45097
$ws.Range("G8").Value = 5687.73

# Row 9
$ws.Range("A9").Value = 42664
$ws.Range("B9").Value = "Giovanna Fogaça"
$ws.Range("C9").Value = "Recursos Humanos"
$ws.Range("D9").Value = "Problemas pessoais"
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 45083
$ws.Range("G9").Value = 4117.19

# Row 10
$ws.Range("A10").Value = 65567
$ws.Range("B10").Value = "Luiz Fernando da Conceição"
$ws.Range("C10").Value = "Jurídico"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 45106
$ws.Range("G10").Value = 3592.41

# Row 11
$ws.Range("A11").Value = 73384
$ws.Range("B11").Value = "Thomas Silva"
$ws.Range("C11").Value = "Recursos Humanos"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 45103
$ws.Range("G11").Value = 4792.7
